# Applies the swap/update of artfynd rows 16-29 (taxon/coordinate re-pairing,
# substrate cell add/remove, and taxon sort-order refresh) per the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16
$ws.Range("A16").Value = 130827881
$ws.Range("B16").Value = 8440
$ws.Range("E16").Value = 106554
$ws.Range("F16").Value = 'Björksplintborre'
$ws.Range("G16").Value = 'Scolytus ratzeburgii'
$ws.Range("H16").Value = 'Janson, 1856'
$ws.Range("Q16").Value = 344495
$ws.Range("R16").Value = 6433268
$ws.Range("AJ16").ClearContents()
$ws.Range("AK16").ClearContents()
$ws.Range("AO16").ClearContents()

# Row 17
$ws.Range("A17").Value = 130827869
$ws.Range("B17").Value = 75221
$ws.Range("E17").Value = 6428
$ws.Range("F17").Value = 'Rostfläck'
$ws.Range("G17").Value = 'Arthonia vinosa'
$ws.Range("H17").Value = 'Leight.'
$ws.Range("Q17").Value = 344518
$ws.Range("R17").Value = 6433262
$ws.Range("AJ17").Value = 'gran'
$ws.Range("AK17").Value = 'Picea abies'
$ws.Range("AO17").Value = 'Picea abies'

# Row 19
$ws.Range("A19").Value = 130827874
$ws.Range("B19").Value = 75349
$ws.Range("E19").Value = 6426
$ws.Range("F19").Value = 'Kattfotslav'
$ws.Range("G19").Value = 'Felipes leucopellaeus'
$ws.Range("H19").Value = '(Ach.) Frisch & G.Thor'
$ws.Range("Q19").Value = 344473
$ws.Range("R19").Value = 6433342
$ws.Range("AJ19").Value = 'gran'
$ws.Range("AK19").Value = 'Picea abies'
$ws.Range("AO19").Value = 'Picea abies'

# Row 20
$ws.Range("A20").Value = 130827872
$ws.Range("B20").Value = 58256
$ws.Range("E20").Value = 103015
$ws.Range("F20").Value = 'Kungsfågel'
$ws.Range("G20").Value = 'Regulus regulus'
$ws.Range("H20").Value = '(Linnaeus, 1758)'
$ws.Range("Q20").Value = 344461
$ws.Range("R20").Value = 6433331
$ws.Range("AJ20").ClearContents()
$ws.Range("AK20").ClearContents()
$ws.Range("AO20").ClearContents()

# Row 21
$ws.Range("B21").Value = 97628

# Row 22
$ws.Range("B22").Value = 83206

# Row 23
$ws.Range("B23").Value = 75349

# Row 24
$ws.Range("A24").Value = 130827870
$ws.Range("B24").Value = 5197
$ws.Range("E24").Value = 105930
$ws.Range("F24").Value = 'Vågbandad barkbock'
$ws.Range("G24").Value = 'Semanotus undatus'
$ws.Range("H24").Value = '(Linnaeus, 1758)'
$ws.Range("Q24").Value = 344458
$ws.Range("R24").Value = 6433350
$ws.Range("AJ24").Value = 'gran'
$ws.Range("AK24").Value = 'Picea abies'
$ws.Range("AO24").Value = 'Picea abies'

# Row 25
$ws.Range("A25").Value = 130827882
$ws.Range("B25").Value = 8440
$ws.Range("E25").Value = 106554
$ws.Range("F25").Value = 'Björksplintborre'
$ws.Range("G25").Value = 'Scolytus ratzeburgii'
$ws.Range("H25").Value = 'Janson, 1856'
$ws.Range("Q25").Value = 344503
$ws.Range("R25").Value = 6433291
$ws.Range("AJ25").ClearContents()
$ws.Range("AK25").ClearContents()
$ws.Range("AO25").ClearContents()

# Row 26
$ws.Range("A26").Value = 130827873
$ws.Range("B26").Value = 83208
$ws.Range("E26").Value = 306
$ws.Range("F26").Value = 'Kornig nållav'
$ws.Range("G26").Value = 'Chaenotheca chlorella'
$ws.Range("H26").Value = '(Ach.) Müll.Arg.'
$ws.Range("Q26").Value = 344451
$ws.Range("R26").Value = 6433334

# Row 27
$ws.Range("A27").Value = 130827876
$ws.Range("B27").Value = 75349
$ws.Range("E27").Value = 6426
$ws.Range("F27").Value = 'Kattfotslav'
$ws.Range("G27").Value = 'Felipes leucopellaeus'
$ws.Range("H27").Value = '(Ach.) Frisch & G.Thor'
$ws.Range("Q27").Value = 344449
$ws.Range("R27").Value = 6433318

# Row 28
$ws.Range("B28").Value = 57881

# Row 29
$ws.Range("B29").Value = 58043
